# Update workbook for "Add data for 2022-11-05" commit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also update the tab name) from 2022-10-27 to 2022-10-28
$ws.Name = "Through 2022-10-28"

# Update the October label cell to reflect new "through" date
$ws.Range("A11").Value = "October (through 10-28)"

# Update October row (row 11) values for each year column (B..I)
$ws.Range("B11").Value = 28
$ws.Range("C11").Value = 47
$ws.Range("D11").Value = 75
$ws.Range("E11").Value = 59
$ws.Range("F11").Value = 56
$ws.Range("G11").Value = 136
$ws.Range("H11").Value = 172
$ws.Range("I11").Value = 101

# Update Total row (row 12) values for each year column (B..I)
$ws.Range("B12").Value = 254
$ws.Range("C12").Value = 476
$ws.Range("D12").Value = 702
$ws.Range("E12").Value = 607
$ws.Range("F12").Value = 478
$ws.Range("G12").Value = 1037
$ws.Range("H12").Value = 1419
$ws.Range("I12").Value = 1377
